$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.456.99'
$ws.Range('E2').Value = '  -2.22%  '
$ws.Range('D3').Value = '1.836.33'
$ws.Range('E3').Value = '  -2.89%  '
$ws.Range('E4').Value = '  -1.02%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '332.09'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -1.36%  '
$ws.Range('E6').Value = '  -1.04%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.4607'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  -3.27%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3813'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -3.85%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '46.37'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -1.86%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.07872'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -2.25%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.9735'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -4.84%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '21.09'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = '  -4.11%  '
$ws.Range('D13').Value = '1.846.49'
$ws.Range('E13').Value = '  -1.96%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '5.890'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  -2.60%  '
$ws.Range('E15').Value = '  -3.14%  '
$ws.Range('E16').Value = '  -1.19%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '87.82'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  -0.92%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.06616'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -2.23%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.00001029'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -2.22%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '16.88'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -1.47%  '
$ws.Range('E21').Value = '  -1.03%  '
$ws.Range('D22').Value = '27.459.33'
$ws.Range('E22').Value = '  -2.14%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '5.333'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -3.75%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '10.84'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -1.93%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.298'
$ws.Range('D25').Style = "Normal"
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '157.22'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -2.17%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '19.32'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -3.63%  '
$ws.Range('E28').Value = '  -2.32%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '5.333'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -3.96%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '118.60'
$ws.Range('D30').Style = "Normal"
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.9530'
$ws.Range('D31').Style = "Normal"
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.09267'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -3.62%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '3.570'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  -1.88%  '
$ws.Range('E34').Value = '  -2.79%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.319'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  -3.92%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.05933'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -2.57%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.02183'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  -3.43%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '8.064'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  -2.06%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '1.156'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -4.09%  '
$ws.Range('E40').Value = '  -3.04%  '
$ws.Range('E41').Value = '  -3.27%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.237'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -2.10%  '
$ws.Range('B44').Value = 'Decentraland'
$ws.Range('C44').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.5481'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -3.48%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '11.98'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  -2.38%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '1.864'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -3.90%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.06649'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  -2.64%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '109.78'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -2.38%  '
$ws.Range('E49').Value = '  -3.19%  '
$ws.Range('B50').Value = 'BabyDogeCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.00000000288'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  -4.76%  '
$ws.Range('B51').Value = 'PaxDollar'
$ws.Range('C51').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.001'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -1.09%  '
